$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 3

# Add new rows 4 and 5 with same formatting as row 3 (style index 1 on column A)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 1

# Copy formatting from A3 (bold, bordered, centered) to A4 and A5
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats = -4122

$excel.CutCopyMode = 0
